$wb = $excel.ActiveWorkbook

# Update the "主键" (primary key) description to "主键，自动递增" (primary key, auto-increment)
# on both the production unit and the receiving unit storage sheets.
$wsProduction = $wb.Worksheets.Item("生产单位仓储表")
$wsProduction.Range("G2").Value = "主键，自动递增"
$wsProduction.Range("G2").Select()

$wsReceiving = $wb.Worksheets.Item("接受单位仓储表")
$wsReceiving.Range("G2").Value = "主键，自动递增"
$wsReceiving.Range("G2").Select()
